$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry text (becomes a new shared string on save)
$newComment = "pientä css tuunausta, suuremmat fontit jne, Dataview komponenttiin redux storen resetointi componentWillUnmount avulla"

# --- Move the "tunnit yht." SUM row from row 75 down to row 86 ---
# Copy formatting of the old totals row (A75:B75) onto the new location (A86:B86)
$ws.Range("A75:B75").Copy()
$ws.Range("A86:B86").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(86).RowHeight = 14.25

$ws.Cells.Item(86, 1).Value = "tunnit yht."
$ws.Cells.Item(86, 2).Formula = "=SUM(B2:B75)"

# --- Turn the old row 75 (previously the totals row) into a normal log entry ---
# Copy formatting from a regular dated log row (row 73) onto row 75
$ws.Range("A73:D73").Copy()
$ws.Range("A75:D75").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(75).AutoFit() | Out-Null

$ws.Cells.Item(75, 1).Value = 44575
$ws.Cells.Item(75, 2).Value = 1
$ws.Cells.Item(75, 3).Value = $newComment
$ws.Cells.Item(75, 4).Value = "client"

# --- Update the visible selection/scroll position ---
$ws.Range("C76").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 65

$excel.CutCopyMode = $false
